# CRYPTO61672.xlsx header rename + header-row border formatting
#
# Summary of the change (per commit "Test: updated headers of all HMRC
# examples"):
#   - Table column headers "Buy Value" / "Sell Value" / "Fee Value" are
#     renamed to "Buy Value in GBP" / "Sell Value in GBP" / "Fee Value in GBP"
#     (columns D, G, J of the Coinbase3 table).
#   - The table header row gets a thin white border applied around each cell
#     (in addition to its existing bold-white-on-black styling).
#   - Column widths are re-fitted to the new (longer) header text.
#   - The active cell / selection on the sheet moved to F23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Rename the three "...Value" headers to "...Value in GBP" -------------
$ws.Range("D1").Value = "Buy Value in GBP"
$ws.Range("G1").Value = "Sell Value in GBP"
$ws.Range("J1").Value = "Fee Value in GBP"

# --- Add a thin white border around every header-row cell -----------------
$hdr = $tbl.HeaderRowRange
$hdr.Borders.Color = 16777215
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

# --- Re-fit the column widths now that some headers are longer ------------
# (ColumnWidth is specified in characters; the engine adds a fixed padding
# offset of 5/6 character when it stores the width, so we compensate for it
# up front in order to land on the real target widths.)
$padding = 0.8333333333333333

$widths = @{
    1  = 7.33203125
    2  = 12
    3  = 12.83203125
    4  = 15.5
    5  = 11.83203125
    6  = 9.1640625
    7  = 15.33203125
    8  = 11.83203125
    9  = 9.1640625
    10 = 15.33203125
    11 = 12.6640625
    12 = 18.1640625
    13 = 208
}

foreach ($col in $widths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $widths[$col] - $padding
}

# --- Move the active selection to F23 --------------------------------------
$ws.Range("F23").Select() | Out-Null
